# Updated symbol list on Fri Dec 30 07:48:53 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as plain TEXT
# (inline strings in the source workbook), not as real numbers. Assigning
# a numeric-looking string straight to Range.Value makes Excel coerce it
# to a Number cell, which would change the cell's type. To keep these as
# Text cells we prefix the literal with a single quote (forces text entry,
# same as typing '244.19 into Excel) and then clear the resulting cell
# format (ClearFormats) so the quote-prefix flag doesn't leave a stray
# cell-style behind - the cell ends up plain Text with no style, matching
# the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $text
    $rng.ClearFormats()
}

# --- Price (column D) updates -------------------------------------------
Set-TextValue "D2"  "244.19"
Set-TextValue "D4"  "5.202"
Set-TextValue "D6"  "6.485"
Set-TextValue "D7"  "3.128"
Set-TextValue "D8"  "0.8119"
Set-TextValue "D9"  "0.8440"
Set-TextValue "D10" "0.1350"
Set-TextValue "D11" "0.06942"
Set-TextValue "D12" "0.03138"
Set-TextValue "D13" "0.02877"
Set-TextValue "D14" "0.09376"
Set-TextValue "D15" "3.771"
Set-TextValue "D16" "0.001509"
Set-TextValue "D17" "0.04685"
Set-TextValue "D18" "0.0005999"
Set-TextValue "D19" "0.006280"
Set-TextValue "D21" "0.004279"
Set-TextValue "D22" "0.00008715"
Set-TextValue "D25" "0.3170"
Set-TextValue "D26" "0.1339"
Set-TextValue "D27" "0.1361"

# --- Row 18: "Worst in 24h" tag moved onto OneONE ------------------------
$ws.Range("E18").Value = "17OneONEWorstin24h"

# --- Rows 41-43: rankings rotated (Kick -> 41, BKEX -> 42, CEJI -> 43) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006321"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1049"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002926"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"

# --- Remaining price updates ----------------------------------------------
Set-TextValue "D44" "0.007388"
Set-TextValue "D45" "0.00005279"
Set-TextValue "D49" "0.00002103"
